# AWS-AI-Challenge-Brogrammers.pptx — slide 10 update
#
# 1. Retitle the slide: "IBM-Recommender" -> "AWS-Recommender"
# 2. Drop the small URL-bar screenshot ("Picture 6", the old
#    ibm-recommender.herokuapp.com address bar capture).
# 3. The big diagram picture ("Picture 9") is removed too, and the
#    surviving picture on the slide ends up showing what used to be
#    "Picture 6"'s image, resized/renamed/re-ordered into "Picture 9"'s
#    old slot -- i.e. the net effect is: one picture shape remains,
#    in the same position/size "Picture 9" used to occupy, but
#    displaying the image that was embedded for "Picture 6".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- 1. Title text -----------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "AWS-Recommender"

# --- 2 & 3. Pictures -----------------------------------------------------
$urlBarPic = $s.Shapes.Item("Picture 6")
$diagramPic = $s.Shapes.Item("Picture 9")

# Remember where/how big the diagram picture was before removing it.
$left = $diagramPic.Left
$top = $diagramPic.Top
$width = $diagramPic.Width
$height = $diagramPic.Height
$name = $diagramPic.Name

# Remove the old diagram picture entirely.
$diagramPic.Delete()

# Re-purpose the remaining picture: move/resize it into the slot the
# diagram picture used to occupy, rename it to match, and restack it so
# it is the last shape on the slide (matching the diagram picture's old
# z-order position).
$urlBarPic.Left = $left
$urlBarPic.Top = $top
$urlBarPic.Width = $width
$urlBarPic.Height = $height
$urlBarPic.Name = $name
$urlBarPic.ZOrder(2)  # msoBringForward
